$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 112083111
$ws.Range("B7").Value = 90814
$ws.Range("E7").Value = 4364
$ws.Range("F7").Value = "Dropptaggsvamp"
$ws.Range("G7").Value = "Hydnellum ferrugineum"
$ws.Range("H7").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q7").Value = 412205
$ws.Range("R7").Value = 6655989

# Row 8
$ws.Range("B8").Value = 89517

# Row 9
$ws.Range("A9").Value = 112083127
$ws.Range("B9").Value = 77738
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 6450
$ws.Range("F9").Value = "Skuggblåslav"
$ws.Range("G9").Value = "Hypogymnia vittata"
$ws.Range("H9").Value = "(Ach.) Parrique"
$ws.Range("Q9").Value = 413052
$ws.Range("R9").Value = 6656343

# Row 10
$ws.Range("B10").Value = 77321

# Row 11
$ws.Range("A11").Value = 112083110
$ws.Range("B11").Value = 78242
$ws.Range("E11").Value = 6453
$ws.Range("F11").Value = "Vedskivlav"
$ws.Range("G11").Value = "Hertelidea botryosa"
$ws.Range("H11").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q11").Value = 412206
$ws.Range("R11").Value = 6656051

# Row 12
$ws.Range("A12").Value = 112083112
$ws.Range("B12").Value = 79580
$ws.Range("E12").Value = 1049
$ws.Range("F12").Value = "Kortskaftad ärgspik"
$ws.Range("G12").Value = "Microcalicium ahlneri"
$ws.Range("H12").Value = "Tibell"
$ws.Range("Q12").Value = 412284
$ws.Range("R12").Value = 6656072

# Row 13
$ws.Range("A13").Value = 112083118
$ws.Range("B13").Value = 94301
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 53
$ws.Range("F13").Value = "Vedtrappmossa"
$ws.Range("G13").Value = "Crossocalyx hellerianus"
$ws.Range("H13").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q13").Value = 412577
$ws.Range("R13").Value = 6656304

# Row 14
$ws.Range("A14").Value = 112083126
$ws.Range("B14").Value = 78671
$ws.Range("E14").Value = 229497
$ws.Range("F14").Value = "Korallblylav"
$ws.Range("G14").Value = "Parmeliella triptophylla"
$ws.Range("H14").Value = "(Ach.) Müll.Arg."
$ws.Range("Q14").Value = 413017
$ws.Range("R14").Value = 6656342
